# Update "paises" (countries) COVID-19 data workbook:
#  - refresh the "last updated" timestamp banner
#  - update country case figures; because the table is kept sorted in
#    descending order by "Casos totales", updated figures make a handful
#    of countries overtake their neighbours, so the country label shown
#    in a given row also changes for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 16:50"

# Each row is written as: country label (col A) then the 7 figures for
# columns B..H (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes).

# Chile overtakes Dinamarca
$row = 28
$vals = @("Chile", 3737, 333, 427, 3288, 31, 4, 22)
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item($row, $i + 1).Value = $vals[$i] }

$row = 29
$vals = @("Dinamarca", 3672, 286, 1193, 2340, 153, 16, 139)
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item($row, $i + 1).Value = $vals[$i] }

# Republica Dominicana overtakes Serbia, Panama, Sudafrica, Peru
$row = 45
$vals = @("Republica Dominicana", 1488, 108, 16, 1404, 147, 8, 68)
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item($row, $i + 1).Value = $vals[$i] }

$row = 46
$vals = @("Serbia", 1476, 305, 42, 1395, 81, 8, 39)
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item($row, $i + 1).Value = $vals[$i] }

$row = 47
$vals = @("Panama", 1475, 0, 9, 1429, 50, 0, 37)
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item($row, $i + 1).Value = $vals[$i] }

$row = 48
$vals = @("Sudafrica", 1462, 0, 95, 1362, 7, 0, 5)
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item($row, $i + 1).Value = $vals[$i] }

$row = 49
$vals = @("Peru", 1414, 0, 537, 822, 51, 0, 55)
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item($row, $i + 1).Value = $vals[$i] }

# Bosnia y Herzegovina keeps its rank, figures updated
$row = 71
$vals = @("Bosnia y Herzegovina", 574, 41, 27, 530, 4, 1, 17)
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item($row, $i + 1).Value = $vals[$i] }

# Madagascar overtakes Banglades, Aruba, Monaco
$row = 128
$vals = @("Madagascar", 65, 6, 0, 65, 6, 0, 0)
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item($row, $i + 1).Value = $vals[$i] }

$row = 129
$vals = @("Banglades", 61, 5, 26, 29, 1, 0, 6)
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item($row, $i + 1).Value = $vals[$i] }

$row = 130
$vals = @("Aruba", 60, 0, 1, 59, 0, 0, 0)
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item($row, $i + 1).Value = $vals[$i] }

$row = 131
$vals = @("Monaco", 60, 0, 2, 57, 2, 0, 1)
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item($row, $i + 1).Value = $vals[$i] }

# Liberia keeps its rank, figures updated
$row = 190
$vals = @("Liberia", 7, 1, 0, 7, 0, 0, 0)
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item($row, $i + 1).Value = $vals[$i] }

# Belice overtakes Botsuana, Gambia
$row = 199
$vals = @("Belice", 4, 1, 0, 4, 0, 0, 0)
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item($row, $i + 1).Value = $vals[$i] }

$row = 200
$vals = @("Botsuana", 4, 0, 0, 3, 0, 0, 1)
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item($row, $i + 1).Value = $vals[$i] }

$row = 201
$vals = @("Gambia", 4, 0, 2, 1, 0, 0, 1)
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item($row, $i + 1).Value = $vals[$i] }
